$wb = $excel.ActiveWorkbook

# --- Productdata sheet: update AverageDemand (column G) for rows 2-5 ---
$wsProd = $wb.Worksheets.Item("Productdata")
$wsProd.Range("G2").Value = 49
$wsProd.Range("G3").Value = 21
$wsProd.Range("G4").Value = 35
$wsProd.Range("G5").Value = 70

# Re-assert the (unchanged) blank StandardDevDemands column so the
# round-trip keeps these cells empty instead of drifting to a stray value.
$wsProd.Range("H2:H11").Value = ""

# --- ForecastedAverageDemand sheet: update rows 9-11, columns B-E ---
$wsFAD = $wb.Worksheets.Item("ForecastedAverageDemand")
$wsFAD.Range("B9").Value = 70
$wsFAD.Range("C9").Value = 30
$wsFAD.Range("D9").Value = 50
$wsFAD.Range("E9").Value = 100

$wsFAD.Range("B10").Value = 70
$wsFAD.Range("C10").Value = 30
$wsFAD.Range("D10").Value = 50
$wsFAD.Range("E10").Value = 100

$wsFAD.Range("B11").Value = 70
$wsFAD.Range("C11").Value = 30
$wsFAD.Range("D11").Value = 50
$wsFAD.Range("E11").Value = 100

# --- ForcastedStandardDeviation sheet: update rows 9-11, columns B-E ---
$wsFSD = $wb.Worksheets.Item("ForcastedStandardDeviation")
$wsFSD.Range("B9").Value = 7.166424999999998
$wsFSD.Range("C9").Value = 3.071324999999999
$wsFSD.Range("D9").Value = 5.118874999999999
$wsFSD.Range("E9").Value = 10.23775

$wsFSD.Range("B10").Value = 8.1997825
$wsFSD.Range("C10").Value = 3.5141925
$wsFSD.Range("D10").Value = 5.856987499999999
$wsFSD.Range("E10").Value = 11.713975

$wsFSD.Range("B11").Value = 9.129804249999998
$wsFSD.Range("C11").Value = 3.912773249999999
$wsFSD.Range("D11").Value = 6.521288749999998
$wsFSD.Range("E11").Value = 13.0425775
